$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 155.8700226666667
$ws.Range("H2").Value = 467.610068
$ws.Range("I2").Value = 0.4627663557222626
$ws.Range("J2").Value = 0.4864916976605717
$ws.Range("M2").Value = 3.031236
$ws.Range("N2").Value = 9.093708000000001
$ws.Range("O2").Value = 0.6569357730026921
$ws.Range("P2").Value = 0.6780946256479073
$ws.Range("Q2").Value = 472.478824028016
$ws.Range("R2").Value = 4252.309416252145
$ws.Range("S2").Value = 0.3040077736160434
$ws.Range("T2").Value = 0.3298874056059602
$ws.Range("G3").Value = 155.8700226666667
$ws.Range("H3").Value = 467.610068
$ws.Range("I3").Value = 0.4627663557222626
$ws.Range("J3").Value = 0.4864916976605717
$ws.Range("O3").Value = 0.2314409052885859
$ws.Range("P3").Value = 0.2388952474211406
$ws.Range("Q3").Value = 166.4560391694218
$ws.Range("R3").Value = 1498.104352524796
$ws.Range("S3").Value = 0.1071030643054602
$ws.Range("T3").Value = 0.116220554480953
$ws.Range("G4").Value = 155.8700226666667
$ws.Range("H4").Value = 467.610068
$ws.Range("I4").Value = 0.4627663557222626
$ws.Range("J4").Value = 0.4864916976605717
$ws.Range("M4").Value = 0.04253966666666667
$ws.Range("N4").Value = 0.127619
$ws.Range("O4").Value = 0.009219285072143351
$ws.Range("P4").Value = 0.009516223528461688
$ws.Range("Q4").Value = 6.630658807565778
$ws.Range("R4").Value = 59.67592926809201
$ws.Range("S4").Value = 0.004266374955200435
$ws.Range("T4").Value = 0.004629563739678802
$ws.Range("G5").Value = 155.8700226666667
$ws.Range("H5").Value = 467.610068
$ws.Range("I5").Value = 0.4627663557222626
$ws.Range("J5").Value = 0.4864916976605717
$ws.Range("M5").Value = 0.4319364999999999
$ws.Range("N5").Value = 0.8638729999999999
$ws.Range("O5").Value = 0.09361017700884301
$ws.Range("P5").Value = 0.06441680759293508
$ws.Range("Q5").Value = 67.32595204556065
$ws.Range("R5").Value = 403.955712273364
$ws.Range("S5").Value = 0.04331964047289821
$ws.Range("T5").Value = 0.03133824208376139
$ws.Range("G6").Value = 155.8700226666667
$ws.Range("H6").Value = 467.610068
$ws.Range("I6").Value = 0.4627663557222626
$ws.Range("J6").Value = 0.4864916976605717
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.04057666666666667
$ws.Range("N6").Value = 0.12173
$ws.Range("O6").Value = 0.008793859627735762
$ws.Range("P6").Value = 0.009077095809555328
$ws.Range("Q6").Value = 6.324685953071111
$ws.Range("R6").Value = 56.92217357764
$ws.Range("S6").Value = 0.004069502372660411
$ws.Range("T6").Value = 0.004415931750218233
$ws.Range("I7").Value = 0.3897411505765819
$ws.Range("J7").Value = 0.4097225989911443
$ws.Range("M7").Value = 3.031236
$ws.Range("N7").Value = 9.093708000000001
$ws.Range("O7").Value = 0.6569357730026921
$ws.Range("P7").Value = 0.6780946256479073
$ws.Range("Q7").Value = 397.920977233416
$ws.Range("R7").Value = 3581.288795100744
$ws.Range("S7").Value = 0.2560349040249855
$ws.Range("T7").Value = 0.2778306923823876
$ws.Range("I8").Value = 0.3897411505765819
$ws.Range("J8").Value = 0.4097225989911443
$ws.Range("O8").Value = 0.2314409052885859
$ws.Range("P8").Value = 0.2388952474211406
$ws.Range("S8").Value = 0.09020204471765919
$ws.Range("T8").Value = 0.09788078166002219
$ws.Range("I9").Value = 0.3897411505765819
$ws.Range("J9").Value = 0.4097225989911443
$ws.Range("M9").Value = 0.04253966666666667
$ws.Range("N9").Value = 0.127619
$ws.Range("O9").Value = 0.009219285072143351
$ws.Range("P9").Value = 0.009516223528461688
$ws.Range("Q9").Value = 5.584331187404667
$ws.Range("R9").Value = 50.258980686642
$ws.Range("S9").Value = 0.003593134771510656
$ws.Range("T9").Value = 0.003899011836662001
$ws.Range("I10").Value = 0.3897411505765819
$ws.Range("J10").Value = 0.4097225989911443
$ws.Range("M10").Value = 0.4319364999999999
$ws.Range("N10").Value = 0.8638729999999999
$ws.Range("O10").Value = 0.09361017700884301
$ws.Range("P10").Value = 0.06441680759293508
$ws.Range("Q10").Value = 56.70181872436899
$ws.Range("R10").Value = 340.2109123462139
$ws.Range("S10").Value = 0.03648373809310397
$ws.Range("T10").Value = 0.02639302182568984
$ws.Range("I11").Value = 0.3897411505765819
$ws.Range("J11").Value = 0.4097225989911443
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.04057666666666667
$ws.Range("N11").Value = 0.12173
$ws.Range("O11").Value = 0.008793859627735762
$ws.Range("P11").Value = 0.009077095809555328
$ws.Range("Q11").Value = 5.326641295126668
$ws.Range("R11").Value = 47.93977165614
$ws.Range("S11").Value = 0.003427328969322688
$ws.Range("T11").Value = 0.003719091286382634
$ws.Range("G12").Value = 0.2461213333333333
$ws.Range("H12").Value = 0.738364
$ws.Range("I12").Value = 0.0007307156985262189
$ws.Range("J12").Value = 0.0007681784042585035
$ws.Range("M12").Value = 3.031236
$ws.Range("N12").Value = 9.093708000000001
$ws.Range("O12").Value = 0.6569357730026921
$ws.Range("P12").Value = 0.6780946256479073
$ws.Range("Q12").Value = 0.746051845968
$ws.Range("R12").Value = 6.714466613712001
$ws.Range("S12").Value = 0.0004800332822565238
$ws.Range("T12").Value = 0.0005208976474664768
$ws.Range("G13").Value = 0.2461213333333333
$ws.Range("H13").Value = 0.738364
$ws.Range("I13").Value = 0.0007307156985262189
$ws.Range("J13").Value = 0.0007681784042585035
$ws.Range("O13").Value = 0.2314409052885859
$ws.Range("P13").Value = 0.2388952474211406
$ws.Range("Q13").Value = 0.2628368277675556
$ws.Range("R13").Value = 2.365531449908
$ws.Range("S13").Value = 0.0001691175027754895
$ws.Range("T13").Value = 0.0001835141699489122
$ws.Range("G14").Value = 0.2461213333333333
$ws.Range("H14").Value = 0.738364
$ws.Range("I14").Value = 0.0007307156985262189
$ws.Range("J14").Value = 0.0007681784042585035
$ws.Range("M14").Value = 0.04253966666666667
$ws.Range("N14").Value = 0.127619
$ws.Range("O14").Value = 0.009219285072143351
$ws.Range("P14").Value = 0.009516223528461688
$ws.Range("Q14").Value = 0.01046991947955556
$ws.Range("R14").Value = 0.094229275316
$ws.Range("S14").Value = 0.000006736676331403571
$ws.Range("T14").Value = 0.000007310157404660925
$ws.Range("G15").Value = 0.2461213333333333
$ws.Range("H15").Value = 0.738364
$ws.Range("I15").Value = 0.0007307156985262189
$ws.Range("J15").Value = 0.0007681784042585035
$ws.Range("M15").Value = 0.4319364999999999
$ws.Range("N15").Value = 0.8638729999999999
$ws.Range("O15").Value = 0.09361017700884301
$ws.Range("P15").Value = 0.06441680759293508
$ws.Range("Q15").Value = 0.1063087872953333
$ws.Range("R15").Value = 0.637852723772
$ws.Range("S15").Value = 0.00006840242588217971
$ws.Range("T15").Value = 0.00004948360046416792
$ws.Range("G16").Value = 0.2461213333333333
$ws.Range("H16").Value = 0.738364
$ws.Range("I16").Value = 0.0007307156985262189
$ws.Range("J16").Value = 0.0007681784042585035
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04057666666666667
$ws.Range("N16").Value = 0.12173
$ws.Range("O16").Value = 0.008793859627735762
$ws.Range("P16").Value = 0.009077095809555328
$ws.Range("Q16").Value = 0.009986783302222224
$ws.Range("R16").Value = 0.08988104972000001
$ws.Range("S16").Value = 0.000006425811280622453
$ws.Range("T16").Value = 0.000006972828974285761
$ws.Range("G17").Value = 49.2786865
$ws.Range("H17").Value = 98.557373
$ws.Range("I17").Value = 0.1463047080910041
$ws.Range("J17").Value = 0.1025370217386683
$ws.Range("M17").Value = 3.031236
$ws.Range("N17").Value = 9.093708000000001
$ws.Range("O17").Value = 0.6569357730026921
$ws.Range("P17").Value = 0.6780946256479073
$ws.Range("Q17").Value = 149.375328551514
$ws.Range("R17").Value = 896.2519713090841
$ws.Range("S17").Value = 0.096112796503697
$ws.Range("T17").Value = 0.06952980337093365
$ws.Range("G18").Value = 49.2786865
$ws.Range("H18").Value = 98.557373
$ws.Range("I18").Value = 0.1463047080910041
$ws.Range("J18").Value = 0.1025370217386683
$ws.Range("O18").Value = 0.2314409052885859
$ws.Range("P18").Value = 0.2388952474211406
$ws.Range("Q18").Value = 52.62548134610517
$ws.Range("R18").Value = 315.752888076631
$ws.Range("S18").Value = 0.03386089408856428
$ws.Range("T18").Value = 0.02449560717808605
$ws.Range("G19").Value = 49.2786865
$ws.Range("H19").Value = 98.557373
$ws.Range("I19").Value = 0.1463047080910041
$ws.Range("J19").Value = 0.1025370217386683
$ws.Range("M19").Value = 0.04253966666666667
$ws.Range("N19").Value = 0.127619
$ws.Range("O19").Value = 0.009219285072143351
$ws.Range("P19").Value = 0.009516223528461688
$ws.Range("Q19").Value = 2.096298897481167
$ws.Range("R19").Value = 12.577793384887
$ws.Range("S19").Value = 0.001348824811287685
$ws.Range("T19").Value = 0.0009757652188079033
$ws.Range("G20").Value = 49.2786865
$ws.Range("H20").Value = 98.557373
$ws.Range("I20").Value = 0.1463047080910041
$ws.Range("J20").Value = 0.1025370217386683
$ws.Range("M20").Value = 0.4319364999999999
$ws.Range("N20").Value = 0.8638729999999999
$ws.Range("O20").Value = 0.09361017700884301
$ws.Range("P20").Value = 0.06441680759293508
$ws.Range("Q20").Value = 21.28526337140725
$ws.Range("R20").Value = 85.14105348562899
$ws.Range("S20").Value = 0.013695609621626
$ws.Range("T20").Value = 0.0066051076004924
$ws.Range("G21").Value = 49.2786865
$ws.Range("H21").Value = 98.557373
$ws.Range("I21").Value = 0.1463047080910041
$ws.Range("J21").Value = 0.1025370217386683
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.04057666666666667
$ws.Range("N21").Value = 0.12173
$ws.Range("O21").Value = 0.008793859627735762
$ws.Range("P21").Value = 0.009077095809555328
$ws.Range("Q21").Value = 1.999564835881667
$ws.Range("R21").Value = 11.99738901529
$ws.Range("S21").Value = 0.001286583065829147
$ws.Range("T21").Value = 0.0009307383703483499
$ws.Range("G22").Value = 0.1539513333333333
$ws.Range("H22").Value = 0.461854
$ws.Range("I22").Value = 0.0004570699116250634
$ws.Range("J22").Value = 0.0004805032053572586
$ws.Range("M22").Value = 3.031236
$ws.Range("N22").Value = 9.093708000000001
$ws.Range("O22").Value = 0.6569357730026921
$ws.Range("P22").Value = 0.6780946256479073
$ws.Range("Q22").Value = 0.466662823848
$ws.Range("R22").Value = 4.199965414632
$ws.Range("S22").Value = 0.0003002655757096832
$ws.Range("T22").Value = 0.0003258266411593498
$ws.Range("G23").Value = 0.1539513333333333
$ws.Range("H23").Value = 0.461854
$ws.Range("I23").Value = 0.0004570699116250634
$ws.Range("J23").Value = 0.0004805032053572586
$ws.Range("O23").Value = 0.2314409052885859
$ws.Range("P23").Value = 0.2388952474211406
$ws.Range("Q23").Value = 0.1644070407708889
$ws.Range("R23").Value = 1.479663366938
$ws.Range("S23").Value = 0.0001057846741266786
$ws.Range("T23").Value = 0.0001147899321304734
$ws.Range("G24").Value = 0.1539513333333333
$ws.Range("H24").Value = 0.461854
$ws.Range("I24").Value = 0.0004570699116250634
$ws.Range("J24").Value = 0.0004805032053572586
$ws.Range("M24").Value = 0.04253966666666667
$ws.Range("N24").Value = 0.127619
$ws.Range("O24").Value = 0.009219285072143351
$ws.Range("P24").Value = 0.009516223528461688
$ws.Range("Q24").Value = 0.006549038402888889
$ws.Range("R24").Value = 0.058941345626
$ws.Range("S24").Value = 0.000004213857813170828
$ws.Range("T24").Value = 0.000004572575908322002
$ws.Range("G25").Value = 0.1539513333333333
$ws.Range("H25").Value = 0.461854
$ws.Range("I25").Value = 0.0004570699116250634
$ws.Range("J25").Value = 0.0004805032053572586
$ws.Range("M25").Value = 0.4319364999999999
$ws.Range("N25").Value = 0.8638729999999999
$ws.Range("O25").Value = 0.09361017700884301
$ws.Range("P25").Value = 0.06441680759293508
$ws.Range("Q25").Value = 0.06649720009033332
$ws.Range("R25").Value = 0.3989832005419999
$ws.Range("S25").Value = 0.00004278639533263841
$ws.Range("T25").Value = 0.0000309524825272871
$ws.Range("G26").Value = 0.1539513333333333
$ws.Range("H26").Value = 0.461854
$ws.Range("I26").Value = 0.0004570699116250634
$ws.Range("J26").Value = 0.0004805032053572586
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.04057666666666667
$ws.Range("N26").Value = 0.12173
$ws.Range("O26").Value = 0.008793859627735762
$ws.Range("P26").Value = 0.009077095809555328
$ws.Range("Q26").Value = 0.009986783302222224
$ws.Range("R26").Value = 0.08988104972000001
$ws.Range("S26").Value = 0.000006425811280622453
$ws.Range("T26").Value = 0.000006972828974285761
